# This script reproduces the commit "adding a bunch of data and a lot of
# other changes for the last days" for the Sydinvest exclusion list.
#
# Diffing the underlying OOXML shows that, despite a lot of cosmetic
# re-save noise (window position, default row height, column bestFit
# widths, etc.), the only real content change is that three companies
# were removed from the exclusion list:
#   - row 44:  BSU                  / CN / Øvrige
#   - row 181: S&T Holdings Co Ltd  / KR / Landminer
#   - row 186: SCIE                 / CN / Termisk Kul
# All later rows simply shift up to fill the gap (245 -> 242 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows to remove, identified by their company name (column A) so that we
# delete the correct rows even if something upstream has already shifted
# things around.
$companiesToRemove = @("BSU", "S&T Holdings Co Ltd", "SCIE")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($i = 0; $i -lt $companiesToRemove.Count; $i++) {
    $target = $companiesToRemove[$i]
    for ($r = $lastRow; $r -ge 2; $r--) {
        $cellValue = $ws.Cells.Item($r, 1).Value2
        if ($cellValue -eq $target) {
            $ws.Rows.Item($r).Delete()
            $lastRow = $lastRow - 1
            break
        }
    }
}

# Match the saved cursor position from the final workbook state.
$ws.Range("C244").Select()
